# This script updates the '想去人数' (F column) values across the
# '展览', '演出', and '全部类型' worksheets, per the authoritative diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 4691  # was 4683
$ws.Range("F5").Value = 2746  # was 2741
$ws.Range("F10").Value = 744  # was 743
$ws.Range("F12").Value = 217  # was 212
$ws.Range("F13").Value = 396  # was 395
$ws.Range("F24").Value = 143  # was 141
$ws.Range("F26").Value = 505  # was 501
$ws.Range("F27").Value = 1655  # was 1654
$ws.Range("F28").Value = 1473  # was 1460
$ws.Range("F31").Value = 1427  # was 1423
$ws.Range("F32").Value = 2298  # was 2285
$ws.Range("F33").Value = 382  # was 381
$ws.Range("F35").Value = 598  # was 597
$ws.Range("F39").Value = 775  # was 772
$ws.Range("F40").Value = 1460  # was 1456
$ws.Range("F41").Value = 198  # was 197
$ws.Range("F44").Value = 25  # was 22
$ws.Range("F45").Value = 76  # was 75

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 88  # was 87

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 4691  # was 4683
$ws.Range("F4").Value = 2746  # was 2741
$ws.Range("F8").Value = 744  # was 743
$ws.Range("F10").Value = 217  # was 212
$ws.Range("F11").Value = 396  # was 395
$ws.Range("F21").Value = 143  # was 141
$ws.Range("F22").Value = 88  # was 87
$ws.Range("F26").Value = 505  # was 501
$ws.Range("F27").Value = 1655  # was 1654
$ws.Range("F28").Value = 1473  # was 1460
$ws.Range("F33").Value = 2298  # was 2285
$ws.Range("F34").Value = 382  # was 381
$ws.Range("F39").Value = 598  # was 597
$ws.Range("F43").Value = 775  # was 772
$ws.Range("F44").Value = 1460  # was 1456
$ws.Range("F46").Value = 198  # was 197
$ws.Range("F48").Value = 76  # was 75

